$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("F1").Value = "Title"
$ws.Range("G1").Value = "Country"
$ws.Range("H1").Value = "Phone"
$ws.Range("I1").Value = "Email"

# Row 2 - Paul Newman
$ws.Range("G2").Value = "Canada"
$ws.Range("I2").Value = "Email Invalid"

# Row 3 - Neil Lynchehaun
$ws.Range("G3").Value = "Canada"
$ws.Range("I3").Value = "nlynchehaun@flatironcorp.com"

# Row 5 - Octavio Flores
$ws.Range("G5").Value = "Canada"
$ws.Range("I5").Value = "Email Invalid"

# Row 6 - Jarred Gumbleton
$ws.Range("G6").Value = "Canada"

# Row 7 - Alicia Lopez
$ws.Range("G7").Value = "United States"

# Row 8 - Mike M
$ws.Range("G8").Value = "Canada"
$ws.Range("I8").Value = "mmeacher@flatironcorp.com"

# Row 9 - Donald Dow
$ws.Range("G9").Value = "Canada"
$ws.Range("I9").Value = "Email Invalid"

# Row 10 - Jan Kyrstein
$ws.Range("G10").Value = "Canada"
$ws.Range("I10").Value = "Email Invalid"

# Row 12 - Mark Neis
$ws.Range("G12").Value = "Canada"
$ws.Range("I12").Value = "mneis@flatironcorp.com"

# Row 13 - Joel Jacques
$ws.Range("G13").Value = "Canada"

# Row 14 - Fidel Velarde
$ws.Range("G14").Value = "Canada"
$ws.Range("I14").Value = "fvelarde@flatironcorp.com"

# Row 15 - Jenn Hirschman
$ws.Range("G15").Value = "Canada"
$ws.Range("I15").Value = "Email Invalid"

# Row 16 - Pete Walton
$ws.Range("I16").Value = "Email Invalid"

# Row 17 - Frank Mydlinski
$ws.Range("G17").Value = "Canada"
$ws.Range("I17").Value = "Email Invalid"

# Row 18 - Husted Janet
$ws.Range("G18").Value = "Canada"
$ws.Range("I18").Value = "Email Invalid"

# Row 19 - Hodge Garry
$ws.Range("I19").Value = "Email Invalid"
